$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 8061.875
$ws.Range("I18").Value = 9185
$ws.Range("K18").Value = 9185
$ws.Range("M18").Value = -8901
$ws.Range("H28").Value = 5516.091
$ws.Range("I28").Value = 5854.1113
$ws.Range("K28").Value = 5854.1113
$ws.Range("M28").Value = -5369.1113
$ws.Range("H33").Value = 829.6923
$ws.Range("J33").Value = 1436.6666
$ws.Range("L33").Value = 1436.6666
$ws.Range("N33").Value = -1894.6666
$ws.Range("H55").Value = 268.9
$ws.Range("I55").Value = 114.5
$ws.Range("K55").Value = 114.5
$ws.Range("M55").Value = 99.5
$ws.Range("H98").Value = 31335.32
$ws.Range("I98").Value = 51520.637
$ws.Range("J98").Value = 15475.429
$ws.Range("K98").Value = 51520.637
$ws.Range("L98").Value = 15475.429
$ws.Range("M98").Value = -50022.637
$ws.Range("N98").Value = -18471.429
$ws.Range("H112").Value = 2419.25
$ws.Range("J112").Value = 2704.2222
$ws.Range("L112").Value = 8112.6666
$ws.Range("N112").Value = -10328.6666
$ws.Range("H122").Value = 31335.32
$ws.Range("I122").Value = 51520.637
$ws.Range("J122").Value = 15475.429
$ws.Range("K122").Value = 154561.911
$ws.Range("L122").Value = 46426.287
$ws.Range("M122").Value = -152111.911
$ws.Range("N122").Value = -51326.287
$ws.Range("H135").Value = 3330.2727
$ws.Range("I135").Value = 3251.4285
$ws.Range("K135").Value = 29262.8565
$ws.Range("M135").Value = -26727.8565
$ws.Range("H137").Value = 8888.689
$ws.Range("I137").Value = 11037
$ws.Range("K137").Value = 33111
$ws.Range("M137").Value = -30561

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1140.125
$ws.Range("I2").Value = 692.75
$ws.Range("K2").Value = 692.75
$ws.Range("M2").Value = -579.75
$ws.Range("H116").Value = 1140.125
$ws.Range("I116").Value = 692.75
$ws.Range("K116").Value = 692.75
$ws.Range("M116").Value = 1601.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1140.125
$ws.Range("I3").Value = 692.75
$ws.Range("K3").Value = 692.75
$ws.Range("M3").Value = -578.75
$ws.Range("H99").Value = 15101.637
$ws.Range("I99").Value = 18695.4
$ws.Range("K99").Value = 18695.4
$ws.Range("M99").Value = -17197.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13843.077
$ws.Range("I31").Value = 23252.334
$ws.Range("J31").Value = 5778
$ws.Range("K31").Value = 23252.334
$ws.Range("L31").Value = 5778
$ws.Range("M31").Value = -22957.334
$ws.Range("N31").Value = -6368
$ws.Range("H34").Value = 13843.077
$ws.Range("I34").Value = 23252.334
$ws.Range("J34").Value = 5778
$ws.Range("K34").Value = 23252.334
$ws.Range("L34").Value = 5778
$ws.Range("M34").Value = -23050.334
$ws.Range("N34").Value = -6182
$ws.Range("H58").Value = 4031.7222
$ws.Range("I58").Value = 3860.3845
$ws.Range("K58").Value = 3860.3845
$ws.Range("M58").Value = -3657.3845
$ws.Range("H86").Value = 7117.5293
$ws.Range("J86").Value = 7400.778
$ws.Range("L86").Value = 7400.778
$ws.Range("N86").Value = -9646.778
$ws.Range("H89").Value = 7117.5293
$ws.Range("J89").Value = 7400.778
$ws.Range("L89").Value = 37003.89
$ws.Range("N89").Value = -48235.89
$ws.Range("H99").Value = 10420386
$ws.Range("I99").Value = 12502864
$ws.Range("K99").Value = 12502864
$ws.Range("M99").Value = -12501366
$ws.Range("H121").Value = 15755
$ws.Range("J121").Value = 15755
$ws.Range("L121").Value = 15755
$ws.Range("N121").Value = -18375
$ws.Range("H122").Value = 8590.529
$ws.Range("I122").Value = 10625.077
$ws.Range("K122").Value = 31875.231
$ws.Range("M122").Value = -29425.231
$ws.Range("H126").Value = 10420386
$ws.Range("I126").Value = 12502864
$ws.Range("K126").Value = 37508592
$ws.Range("M126").Value = -37506122
$ws.Range("H132").Value = 3466.3333
$ws.Range("I132").Value = 3466.3333
$ws.Range("K132").Value = 10398.9999
$ws.Range("M132").Value = -7868.999899999999
$ws.Range("H134").Value = 3307.0527
$ws.Range("I134").Value = 3559.7856
$ws.Range("K134").Value = 10679.3568
$ws.Range("M134").Value = -8144.356800000001
$ws.Range("H136").Value = 4031.7222
$ws.Range("I136").Value = 3860.3845
$ws.Range("K136").Value = 11581.1535
$ws.Range("M136").Value = -9031.1535

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 207
$ws.Range("I11").Value = 201
$ws.Range("K11").Value = 603
$ws.Range("M11").Value = -463
$ws.Range("H87").Value = 14796.728
$ws.Range("I87").Value = 10205.4
$ws.Range("J87").Value = 18622.834
$ws.Range("K87").Value = 30616.2
$ws.Range("L87").Value = 55868.50199999999
$ws.Range("M87").Value = -29368.2
$ws.Range("N87").Value = -58364.50199999999
$ws.Range("H90").Value = 14796.728
$ws.Range("I90").Value = 10205.4
$ws.Range("J90").Value = 18622.834
$ws.Range("K90").Value = 91848.59999999999
$ws.Range("L90").Value = 167605.506
$ws.Range("M90").Value = -85608.59999999999
$ws.Range("N90").Value = -180085.506

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 16497
$ws.Range("I126").Value = 61000
$ws.Range("J126").Value = 9079.833000000001
$ws.Range("K126").Value = 183000
$ws.Range("L126").Value = 27239.499
$ws.Range("M126").Value = -180530
$ws.Range("N126").Value = -32179.499

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 18245.75
$ws.Range("I40").Value = 24120
$ws.Range("J40").Value = 11680.412
$ws.Range("K40").Value = 24120
$ws.Range("L40").Value = 11680.412
$ws.Range("M40").Value = -23984
$ws.Range("N40").Value = -11952.412
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H48").Value = 24000
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 24000
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 24000
$ws.Range("M48").ClearContents()
$ws.Range("N48").Value = -25322
$ws.Range("H104").Value = 32000
$ws.Range("J104").Value = 32000
$ws.Range("L104").Value = 32000
$ws.Range("N104").Value = -38988
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5088.175
$ws.Range("I122").Value = 2664.3125
$ws.Range("J122").Value = 6704.0835
$ws.Range("K122").Value = 7992.9375
$ws.Range("L122").Value = 20112.2505
$ws.Range("M122").Value = -5542.9375
$ws.Range("N122").Value = -25012.2505
$ws.Range("H126").Value = 26235.941
$ws.Range("I126").Value = 28401.066
$ws.Range("K126").Value = 85203.198
$ws.Range("M126").Value = -82733.198
$ws.Range("H136").Value = 674229.9399999999
$ws.Range("I136").Value = 1030714.6
$ws.Range("J136").Value = 5821.125
$ws.Range("K136").Value = 3092143.8
$ws.Range("L136").Value = 17463.375
$ws.Range("M136").Value = -3089593.8
$ws.Range("N136").Value = -22563.375

